$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert two new rows before the old row 14 ("end repeat") to make room for
# the new "file" and "image" field type rows.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# Match the row height used by the surrounding rows.
$ws.Rows.Item(14).RowHeight = 13.8
$ws.Rows.Item(15).RowHeight = 13.8

# New row 14: file type
$ws.Range("A14").Value = "file"
$ws.Range("B14").Value = "file"
$ws.Range("C14").Value = "File"

# New row 15: image type
$ws.Range("A15").Value = "image"
$ws.Range("B15").Value = "image"
$ws.Range("C15").Value = "Image"

# Update the selection to match the saved cursor position after editing.
$ws.Range("D18").Select()
